$d = $word.ActiveDocument

# 1) Remove the old "_GoBack" bookmark that currently sits after the
#    "Следующим шагом станет подбор доменного имени..." sentence.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Re-create the "_GoBack" bookmark so that it starts right before the
#    "Создать сайт в веб студии W" heading run and ends right after
#    " студии" (i.e. before the trailing " W"). Doing this via a single
#    Bookmarks.Add over that whole span reproduces the exact
#    bookmarkStart / run-split / bookmarkEnd / run layout from the diff.
$headingRange = $d.Content
$headingRange.Find.Execute("Создать сайт в веб") | Out-Null
$headingStart = $headingRange.Start

$afterHeading = $d.Range($headingRange.End, $d.Content.End)
$afterHeading.Find.Execute(" студии") | Out-Null
$bookmarkEndPos = $afterHeading.End

$bookmarkRange = $d.Range($headingStart, $bookmarkEndPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# 3) Turn the trailing empty paragraph into the SEO keywords paragraph,
#    keeping the keywords split across separate runs exactly as authored.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$keywordsXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">создать сайт, сайт в Мариуполе, </w:t></w:r><w:r><w:t>создать сайт</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>в Мариуполе</w:t></w:r><w:r><w:t>, сделать сайт дёшево</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($keywordsXml)
